# "complete move & add & remove"
# Rebuilds the "Map" grid data (object move/add/remove on the map),
# grows the map to 6x8 (MapInfo), and flips the active sheet/selection
# back to Map / Submenu as recorded in the target workbook.

$wb = $excel.ActiveWorkbook

$wsMap     = $wb.Worksheets.Item("Map")
$wsSubmenu = $wb.Worksheets.Item("Submenu")
$wsInfo    = $wb.Worksheets.Item("MapInfo")

# ---------------------------------------------------------------
# Map sheet: rewrite the grid of object placements (A:F, rows 1-9)
# ---------------------------------------------------------------
$mapData = @{
    1 = @(3,3,0,0,0,0)
    2 = @(3,3,0,0,0,0)
    3 = @(0,0,10,10,0,0)
    4 = @(0,0,10,10,0,0)
    5 = @(0,0,19,19,0,0)
    6 = @(0,0,19,19,0,0)
    7 = @(0,0,0,0,0,0)
    8 = @(0,0,0,0,0,0)
}

foreach ($r in 1..8) {
    $rowVals = $mapData[$r]
    for ($c = 1; $c -le 6; $c++) {
        $wsMap.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# New row 9 only spans columns A-D
$wsMap.Cells.Item(9, 1).Value = 8
$wsMap.Cells.Item(9, 2).Value = 8
$wsMap.Cells.Item(9, 3).Value = 1
$wsMap.Cells.Item(9, 4).Value = 0

# ---------------------------------------------------------------
# MapInfo sheet: grow the declared map size to match
# ---------------------------------------------------------------
$wsInfo.Cells.Item(1, 2).Value = 6   # MapSizeX
$wsInfo.Cells.Item(2, 2).Value = 8   # MapSizeY

# ---------------------------------------------------------------
# Selections / active sheet: Submenu is de-activated (selection left
# at F24), Map becomes the active/selected sheet with B5 selected.
# ---------------------------------------------------------------
$wsSubmenu.Range("F24").Select()

$wsMap.Activate()
$wsMap.Range("B5").Select()

Write-Host "Map data rewritten; MapInfo resized to 6x8; Map sheet active."
